$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, pushing the existing rows 284-326 down to 285-327
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row 284 with the new weekly record
$ws.Range("A284").Value = 10
$ws.Range("B284").Value = "Vega Modelo de Temuco"
$ws.Range("C284").Value = "La Araucanía"
$ws.Range("D284").Value = 44951
$ws.Range("E284").Value = 9
$ws.Range("F284").Value = 100112039
$ws.Range("G284").Value = "Ciboulette"
$ws.Range("H284").Value = "Sin especificar"
$ws.Range("I284").Value = "Primera"
$ws.Range("J284").Value = 35
$ws.Range("K284").Value = 5000
$ws.Range("L284").Value = 5000
$ws.Range("M284").Value = 5000
$ws.Range("N284").Value = "$/docena de atados"
$ws.Range("O284").Value = "Provincia de Cautín"
$ws.Range("P284").Value = 1667
$ws.Range("Q284").Value = 3
$ws.Range("R284").Value = "Hortaliza"
